# Auto-generated Excel COM-interop script applying market-data updates
# across multiple worksheets/rows (columns H..N), per commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value2 = 47143
$ws.Range("I16").Value2 = 1409
$ws.Range("J16").Value2 = 70010
$ws.Range("K16").Value2 = 1409
$ws.Range("L16").Value2 = 70010
$ws.Range("M16").Value2 = -1179
$ws.Range("N16").Value2 = -70470

$ws.Range("H40").Value2 = 1005.7
$ws.Range("I40").Value2 = 999
$ws.Range("J40").Value2 = 1008.5714
$ws.Range("K40").Value2 = 999
$ws.Range("L40").Value2 = 1008.5714
$ws.Range("M40").Value2 = -824
$ws.Range("N40").Value2 = -1358.5714

$ws.Range("H81").Value2 = 40000
$ws.Range("I81").Value2 = 30000
$ws.Range("J81").Value2 = 50000
$ws.Range("K81").Value2 = 30000
$ws.Range("L81").Value2 = 50000
$ws.Range("M81").Value2 = -29002
$ws.Range("N81").Value2 = -51996

$ws.Range("H84").Value2 = 40000
$ws.Range("I84").Value2 = 30000
$ws.Range("J84").Value2 = 50000
$ws.Range("K84").Value2 = 90000
$ws.Range("L84").Value2 = 150000
$ws.Range("M84").Value2 = -85008
$ws.Range("N84").Value2 = -159984

$ws.Range("H135").Value2 = 21740898
$ws.Range("I135").Value2 = 283.7
$ws.Range("J135").Value2 = 38464450
$ws.Range("K135").Value2 = 2553.3
$ws.Range("L135").Value2 = 346180050
$ws.Range("M135").Value2 = -18.29999999999973
$ws.Range("N135").Value2 = -346185120

$ws.Range("H137").Value2 = 24788.215
$ws.Range("I137").Value2 = 27219.605
$ws.Range("J137").Value2 = 1690
$ws.Range("K137").Value2 = 81658.815
$ws.Range("L137").Value2 = 5070
$ws.Range("M137").Value2 = -79108.815
$ws.Range("N137").Value2 = -10170

$ws.Range("H141").Value2 = 1782.2106
$ws.Range("I141").Value2 = 1139.7632
$ws.Range("J141").Value2 = 3067.1052
$ws.Range("K141").Value2 = 3419.2896
$ws.Range("L141").Value2 = 9201.3156
$ws.Range("M141").Value2 = 1760.7104
$ws.Range("N141").Value2 = -19561.3156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 771.75
$ws.Range("I61").Value2 = 771.75
$ws.Range("K61").Value2 = 771.75
$ws.Range("M61").Value2 = -559.75

$ws.Range("H97").Value2 = 4124.25
$ws.Range("I97").Value2 = 3000
$ws.Range("J97").Value2 = 5998
$ws.Range("K97").Value2 = 3000
$ws.Range("L97").Value2 = 5998
$ws.Range("M97").Value2 = -2504
$ws.Range("N97").Value2 = -6990

$ws.Range("H102").Value2 = 0
$ws.Range("I102").Value2 = 0
$ws.Range("J102").Value2 = 0
$ws.Range("K102").Value2 = 0
$ws.Range("L102").Value2 = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()

$ws.Range("H110").Value2 = 932.3333
$ws.Range("I110").Value2 = 829.7
$ws.Range("J110").Value2 = 1137.6
$ws.Range("K110").Value2 = 829.7
$ws.Range("L110").Value2 = 1137.6
$ws.Range("M110").Value2 = 1215.3
$ws.Range("N110").Value2 = -5227.6

$ws.Range("H122").Value2 = 1373.1428
$ws.Range("I122").Value2 = 1303
$ws.Range("J122").Value2 = 1466.6666
$ws.Range("K122").Value2 = 3909
$ws.Range("L122").Value2 = 4399.9998
$ws.Range("M122").Value2 = -1459
$ws.Range("N122").Value2 = -9299.9998

$ws.Range("H136").Value2 = 771.75
$ws.Range("I136").Value2 = 771.75
$ws.Range("K136").Value2 = 2315.25
$ws.Range("M136").Value2 = 234.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value2 = 1231.16
$ws.Range("I22").Value2 = 4235.8
$ws.Range("J22").Value2 = 480
$ws.Range("K22").Value2 = 4235.8
$ws.Range("L22").Value2 = 480
$ws.Range("M22").Value2 = -4062.8
$ws.Range("N22").Value2 = -826

$ws.Range("H26").Value2 = 8780.666999999999
$ws.Range("I26").Value2 = 8780.666999999999
$ws.Range("J26").Value2 = 0
$ws.Range("K26").Value2 = 8780.666999999999
$ws.Range("L26").Value2 = 0
$ws.Range("M26").Value2 = -8488.666999999999
$ws.Range("N26").ClearContents()

$ws.Range("H94").Value2 = 3224.4546
$ws.Range("I94").Value2 = 1650
$ws.Range("J94").Value2 = 3814.875
$ws.Range("K94").Value2 = 1650
$ws.Range("L94").Value2 = 3814.875
$ws.Range("M94").Value2 = -1199
$ws.Range("N94").Value2 = -4716.875

$ws.Range("H99").Value2 = 966.6667
$ws.Range("J99").Value2 = 1000
$ws.Range("L99").Value2 = 1000
$ws.Range("N99").Value2 = -3996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 1012.7273
$ws.Range("I16").Value2 = 950.1667
$ws.Range("J16").Value2 = 1087.8
$ws.Range("K16").Value2 = 950.1667
$ws.Range("L16").Value2 = 1087.8
$ws.Range("M16").Value2 = -663.1667
$ws.Range("N16").Value2 = -1661.8

$ws.Range("H31").Value2 = 36309.3
$ws.Range("I31").Value2 = 48627.668
$ws.Range("J31").Value2 = 7566.4443
$ws.Range("K31").Value2 = 48627.668
$ws.Range("L31").Value2 = 7566.4443
$ws.Range("M31").Value2 = -48332.668
$ws.Range("N31").Value2 = -8156.4443

$ws.Range("H33").Value2 = 35365
$ws.Range("I33").Value2 = 9362.5
$ws.Range("J33").Value2 = 70035
$ws.Range("K33").Value2 = 9362.5
$ws.Range("L33").Value2 = 70035
$ws.Range("M33").Value2 = -8983.5
$ws.Range("N33").Value2 = -70793

$ws.Range("H34").Value2 = 36309.3
$ws.Range("I34").Value2 = 48627.668
$ws.Range("J34").Value2 = 7566.4443
$ws.Range("K34").Value2 = 48627.668
$ws.Range("L34").Value2 = 7566.4443
$ws.Range("M34").Value2 = -48425.668
$ws.Range("N34").Value2 = -7970.4443

$ws.Range("H105").Value2 = 1150
$ws.Range("I105").Value2 = 1100
$ws.Range("J105").Value2 = 1200
$ws.Range("K105").Value2 = 1100
$ws.Range("L105").Value2 = 1200
$ws.Range("M105").Value2 = 647
$ws.Range("N105").Value2 = -4694

$ws.Range("H113").Value2 = 1012.7273
$ws.Range("I113").Value2 = 950.1667
$ws.Range("J113").Value2 = 1087.8
$ws.Range("K113").Value2 = 950.1667
$ws.Range("L113").Value2 = 1087.8
$ws.Range("M113").Value2 = 1219.8333
$ws.Range("N113").Value2 = -5427.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value2 = 0
$ws.Range("J139").Value2 = 0
$ws.Range("L139").Value2 = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value2 = 1486.6428
$ws.Range("J93").Value2 = 1804
$ws.Range("L93").Value2 = 1804
$ws.Range("N93").Value2 = -4300

$ws.Range("H122").Value2 = 2265.5557
$ws.Range("I122").Value2 = 1770
$ws.Range("J122").Value2 = 4000
$ws.Range("K122").Value2 = 5310
$ws.Range("L122").Value2 = 12000
$ws.Range("M122").Value2 = -2860
$ws.Range("N122").Value2 = -16900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value2 = 54514.25
$ws.Range("I28").Value2 = 0
$ws.Range("J28").Value2 = 54514.25
$ws.Range("K28").Value2 = 0
$ws.Range("L28").Value2 = 54514.25
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value2 = -55210.25

$ws.Range("H107").Value2 = 349.5
$ws.Range("I107").Value2 = 299.33334
$ws.Range("J107").Value2 = 500
$ws.Range("K107").Value2 = 898.0000200000001
$ws.Range("L107").Value2 = 1500
$ws.Range("M107").Value2 = 1021.99998
$ws.Range("N107").Value2 = -5340
